# luu y khi edit sau khi search thi khong duoc sua ID vi nut edit tim dong de sua tren excel dua vao ID
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Update student data table on sheet 1 ---

# Row 2 (ID 1): Dang / Ngoc Truong Giang -> Pink / Rose, phone/address/picture updated
$ws1.Range("A2").Value = "1"
$ws1.Range("B2").Value = "Pink"
$ws1.Range("C2").Value = "Rose"
$ws1.Range("D2").Value = 36501.94482556713
$ws1.Range("E2").Value = "Female"
$ws1.Range("F2").Value = "1000"
$ws1.Range("G2").Value = "Korea"
$ws1.Range("H2").Value = "D:\Tai lieu mon hoc 2024\Lập trình trực quan\Rose.jpeg"

# Row 3 (ID 2, was ID 5): Tap / Can Binh keeps name, gender/phone/address updated
$ws1.Range("A3").Value = "2"
$ws1.Range("B3").Value = "Tap"
$ws1.Range("C3").Value = "Can Binh"
$ws1.Range("D3").Value = 18239.963888958333
$ws1.Range("E3").Value = "Male"
$ws1.Range("F3").Value = "2000"
$ws1.Range("G3").Value = "Beijing"
$ws1.Range("H3").Value = "D:\Tai lieu mon hoc 2024\Lập trình trực quan\Tap can binh.jpeg"

# Row 4 (ID 3, was Duong Yen Nhi): now Ma Tieu Dao
$ws1.Range("A4").Value = "3"
$ws1.Range("B4").Value = "Ma"
$ws1.Range("C4").Value = "Tieu Dao"
$ws1.Range("D4").Value = 34706.9646059838
$ws1.Range("E4").Value = "Female"
$ws1.Range("F4").Value = "3000"
$ws1.Range("G4").Value = "Quang Dong"
$ws1.Range("H4").Value = "D:\Tai lieu mon hoc 2024\Lập trình trực quan\Ma Tieu Dao.jpeg"

# Row 5 (ID 4, was Vladimir Putin): now Ma Van Thuong
$ws1.Range("A5").Value = "4"
$ws1.Range("B5").Value = "Ma"
$ws1.Range("C5").Value = "Van Thuong"
$ws1.Range("D5").Value = 27599.966736493054
$ws1.Range("E5").Value = "Female"
$ws1.Range("F5").Value = "4000"
$ws1.Range("G5").Value = "Tham Quyen"
$ws1.Range("H5").Value = "D:\Tai lieu mon hoc 2024\Lập trình trực quan\Bi Bi Dong.jpeg"

# --- Update summary stats on sheet 2 ---
$ws2.Range("A2").Value = 4
$ws2.Range("B2").Value = 25
$ws2.Range("C2").Value = 75
